$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.733.07'
$ws.Range('E2').Value = '  +1.03%  '
$ws.Range('D3').Value = '2.517.14'
$ws.Range('E3').Value = '  -0.04%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.35'
$ws.Range('E5').Value = '  +3.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '94.61'
$ws.Range('E6').Value = '  -1.99%  '
$ws.Range('E7').Value = '  -1.30%  '
$ws.Range('E8').Value = '  -0.24%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.531'
$ws.Range('E9').Value = '  -1.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.74'
$ws.Range('E10').Value = '  -1.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0811'
$ws.Range('E11').Value = '  +0.42%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.54'
$ws.Range('E12').Value = '  +0.95%  '
$ws.Range('E13').Value = '  -3.83%  '
$ws.Range('D14').Value = '2.904.12'
$ws.Range('E14').Value = '  -0.15%  '
$ws.Range('D15').Value = '2.528.09'
$ws.Range('E15').Value = '  -1.34%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.21'
$ws.Range('E16').Value = '  +0.73%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.850'
$ws.Range('E17').Value = '  -0.31%  '
$ws.Range('D18').Value = '42.774.32'
$ws.Range('E18').Value = '  +1.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.98'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.71'
$ws.Range('E20').Value = '  +4.78%  '
$ws.Range('D21').Value = '0.0₃0960'
$ws.Range('E21').Value = '  -1.55%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.42'
$ws.Range('E22').Value = '  -2.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '251.00'
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.94'
$ws.Range('E24').Value = '  +1.56%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.02'
$ws.Range('E25').Value = '  +0.10%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.74'
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('E27').Value = '  -0.27%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.42'
$ws.Range('E28').Value = '  +4.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '41.28'
$ws.Range('E29').Value = '  +8.94%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.26'
$ws.Range('E30').Value = '  +0.26%  '
$ws.Range('E31').Value = '  -0.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '155.03'
$ws.Range('E32').Value = '  +0.05%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.30'
$ws.Range('E33').Value = '  +4.98%  '
$ws.Range('E34').Value = '  +2.39%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0784'
$ws.Range('E36').Value = '  -0.38%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.62'
$ws.Range('E37').Value = '  -0.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.112'
$ws.Range('E38').Value = '  -3.98%  '
$ws.Range('E39').Value = '  -0.95%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '23.57'
$ws.Range('E40').Value = '  -2.13%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.32'
$ws.Range('E41').Value = '  +14.54%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0304'
$ws.Range('E42').Value = '  +1.55%  '
$ws.Range('E43').Value = '  +0.20%  '
$ws.Range('E44').Value = '  -2.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.77'
$ws.Range('E45').Value = '  -2.08%  '
$ws.Range('D46').Value = '2.013.03'
$ws.Range('E46').Value = '  -1.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '85.69'
$ws.Range('E47').Value = '  +1.65%  '
$ws.Range('E48').Value = '  -1.04%  '
$ws.Range('D49').Value = '2.757.77'
$ws.Range('E49').Value = '  -0.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '73.49'
$ws.Range('E50').Value = '  +2.04%  '
$ws.Range('E51').Value = '  +1.20%  '
